# Auto-generated: applies the numeric value updates described in the commit diff
# (scheduled-runner refresh of computed market-price / profit columns H-N per leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 9332.48
$ws.Range("J13").Value = 23001.5
$ws.Range("L13").Value = 23001.5
$ws.Range("N13").Value = -23339.5
$ws.Range("H33").Value = 680.625
$ws.Range("I33").Value = 541.6923
$ws.Range("J33").Value = 844.8182
$ws.Range("K33").Value = 541.6923
$ws.Range("L33").Value = 844.8182
$ws.Range("M33").Value = -312.6923
$ws.Range("N33").Value = -1302.8182
$ws.Range("H48").Value = 2999.5
$ws.Range("I48").Value = 2999
$ws.Range("K48").Value = 8997
$ws.Range("M48").Value = -8705
$ws.Range("H56").Value = 2999.5
$ws.Range("I56").Value = 2999
$ws.Range("K56").Value = 8997
$ws.Range("M56").Value = -8463
$ws.Range("H81").Value = 33156.5
$ws.Range("I81").Value = 17298
$ws.Range("J81").Value = 38442.668
$ws.Range("K81").Value = 17298
$ws.Range("L81").Value = 38442.668
$ws.Range("M81").Value = -16300
$ws.Range("N81").Value = -40438.668
$ws.Range("H84").Value = 33156.5
$ws.Range("I84").Value = 17298
$ws.Range("J84").Value = 38442.668
$ws.Range("K84").Value = 51894
$ws.Range("L84").Value = 115328.004
$ws.Range("M84").Value = -46902
$ws.Range("N84").Value = -125312.004
$ws.Range("H116").Value = 3659.258
$ws.Range("I116").Value = 3074.2632
$ws.Range("J116").Value = 4585.5
$ws.Range("K116").Value = 3074.2632
$ws.Range("L116").Value = 4585.5
$ws.Range("M116").Value = 367.7368000000001
$ws.Range("N116").Value = -11469.5
$ws.Range("H132").Value = 3923760.5
$ws.Range("I132").Value = 4349956.5
$ws.Range("J132").Value = 2760
$ws.Range("K132").Value = 13049869.5
$ws.Range("L132").Value = 8280
$ws.Range("M132").Value = -13047339.5
$ws.Range("N132").Value = -13340
$ws.Range("H133").Value = 34655.555
$ws.Range("J133").Value = 34655.555
$ws.Range("L133").Value = 34655.555
$ws.Range("N133").Value = -44775.555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2978.5
$ws.Range("I63").Value = 2767.2222
$ws.Range("J63").Value = 4880
$ws.Range("K63").Value = 2767.2222
$ws.Range("L63").Value = 4880
$ws.Range("M63").Value = -2081.2222
$ws.Range("N63").Value = -6252
$ws.Range("H66").Value = 2978.5
$ws.Range("I66").Value = 2767.2222
$ws.Range("J66").Value = 4880
$ws.Range("K66").Value = 13836.111
$ws.Range("L66").Value = 24400
$ws.Range("M66").Value = -10404.111
$ws.Range("N66").Value = -31264

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3149.0715
$ws.Range("I107").Value = 1935
$ws.Range("J107").Value = 4363.143
$ws.Range("K107").Value = 1935
$ws.Range("L107").Value = 4363.143
$ws.Range("M107").Value = -15
$ws.Range("N107").Value = -8203.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1374.5
$ws.Range("I16").Value = 808.3077
$ws.Range("J16").Value = 2043.6364
$ws.Range("K16").Value = 808.3077
$ws.Range("L16").Value = 2043.6364
$ws.Range("M16").Value = -521.3077
$ws.Range("N16").Value = -2617.6364
$ws.Range("H50").Value = 13497.333
$ws.Range("J50").Value = 13497.333
$ws.Range("L50").Value = 13497.333
$ws.Range("N50").Value = -14747.333
$ws.Range("H58").Value = 8476708
$ws.Range("I58").Value = 1220.8049
$ws.Range("J58").Value = 27781984
$ws.Range("K58").Value = 1220.8049
$ws.Range("L58").Value = 27781984
$ws.Range("M58").Value = -1017.8049
$ws.Range("N58").Value = -27782390
$ws.Range("H105").Value = 4610
$ws.Range("I105").Value = 4338.75
$ws.Range("K105").Value = 4338.75
$ws.Range("M105").Value = -2591.75
$ws.Range("H113").Value = 1374.5
$ws.Range("I113").Value = 808.3077
$ws.Range("J113").Value = 2043.6364
$ws.Range("K113").Value = 808.3077
$ws.Range("L113").Value = 2043.6364
$ws.Range("M113").Value = 1361.6923
$ws.Range("N113").Value = -6383.6364
$ws.Range("H122").Value = 3780.8125
$ws.Range("I122").Value = 3269.5557
$ws.Range("J122").Value = 4438.143
$ws.Range("K122").Value = 9808.667099999999
$ws.Range("L122").Value = 13314.429
$ws.Range("M122").Value = -7358.667099999999
$ws.Range("N122").Value = -18214.429
$ws.Range("H132").Value = 1677.0167
$ws.Range("I132").Value = 1218.0588
$ws.Range("K132").Value = 3654.1764
$ws.Range("M132").Value = -1124.1764
$ws.Range("H134").Value = 1330.711
$ws.Range("I134").Value = 660.375
$ws.Range("K134").Value = 1981.125
$ws.Range("M134").Value = 553.875
$ws.Range("H136").Value = 8476708
$ws.Range("I136").Value = 1220.8049
$ws.Range("J136").Value = 27781984
$ws.Range("K136").Value = 3662.4147
$ws.Range("L136").Value = 83345952
$ws.Range("M136").Value = -1112.4147
$ws.Range("N136").Value = -83351052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 134.71428
$ws.Range("I12").Value = 9.6
$ws.Range("J12").Value = 173.8125
$ws.Range("K12").Value = 28.8
$ws.Range("L12").Value = 521.4375
$ws.Range("M12").Value = 144.2
$ws.Range("N12").Value = -867.4375
$ws.Range("H75").Value = 2257.5557
$ws.Range("I75").Value = 1013
$ws.Range("J75").Value = 2879.8333
$ws.Range("K75").Value = 3039
$ws.Range("L75").Value = 8639.499899999999
$ws.Range("M75").Value = -2041
$ws.Range("N75").Value = -10635.4999
$ws.Range("H78").Value = 2257.5557
$ws.Range("I78").Value = 1013
$ws.Range("J78").Value = 2879.8333
$ws.Range("K78").Value = 9117
$ws.Range("L78").Value = 25918.4997
$ws.Range("M78").Value = -4125
$ws.Range("N78").Value = -35902.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3935.1904
$ws.Range("I122").Value = 2307.375
$ws.Range("J122").Value = 4936.923
$ws.Range("K122").Value = 6922.125
$ws.Range("L122").Value = 14810.769
$ws.Range("M122").Value = -4472.125
$ws.Range("N122").Value = -19710.769
$ws.Range("H126").Value = 2712
$ws.Range("I126").Value = 1881.8182
$ws.Range("J126").Value = 3282.75
$ws.Range("K126").Value = 5645.4546
$ws.Range("L126").Value = 9848.25
$ws.Range("M126").Value = -3175.4546
$ws.Range("N126").Value = -14788.25
$ws.Range("H132").Value = 2672.9434
$ws.Range("I132").Value = 2174.8108
$ws.Range("J132").Value = 3824.875
$ws.Range("K132").Value = 6524.432400000001
$ws.Range("L132").Value = 11474.625
$ws.Range("M132").Value = -3994.432400000001
$ws.Range("N132").Value = -16534.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3114.2856
$ws.Range("H98").Value = 40000
$ws.Range("J98").Value = 40000
$ws.Range("L98").Value = 40000
$ws.Range("N98").Value = -45990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2209460
$ws.Range("I3").Value = 3666766.8
$ws.Range("K3").Value = 3666766.8
$ws.Range("M3").Value = -3666652.8
$ws.Range("H100").Value = 955.55
$ws.Range("J100").Value = 914
$ws.Range("L100").Value = 1828
$ws.Range("N100").Value = -2910
$ws.Range("H136").Value = 1411.9231
$ws.Range("I136").Value = 893.6667
$ws.Range("K136").Value = 2681.0001
$ws.Range("M136").Value = -131.0001000000002
